# Extend the "Akhalkalaki" average-monthly-remuneration table with the
# 2023 column (K), matching the formatting already used by the preceding
# 2022 column (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, borders, alignment, fill, font, ...)
# from the last populated column (J) into the new column (K) before
# writing the new figures, so the new cells render identically to the
# rest of the table.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null

# Year header
$ws.Range("K3").Value = 2023

# Average monthly remuneration of employed persons of business sector
$ws.Range("K4").Value = 716.5

# Women
$ws.Range("K5").Value = 516.70000000000005

# Men
$ws.Range("K6").Value = 793.1
